$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2026-01-25 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-26 Monday", 2) | Out-Null

# Update the division-fact table cells (row, col -> new text)
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "69÷4=17, 1"
$t.Cell(1, 2).Range.Text = "43÷8=5, 3"
$t.Cell(1, 3).Range.Text = "74÷6=12, 2"
$t.Cell(1, 4).Range.Text = "97÷9=10, 7"
$t.Cell(1, 5).Range.Text = "93÷6=15, 3"

$t.Cell(5, 1).Range.Text = "30÷9=3, 3"
$t.Cell(5, 2).Range.Text = "98÷3=32, 2"
$t.Cell(5, 3).Range.Text = "53÷6=8, 5"
$t.Cell(5, 4).Range.Text = "24÷6=4, 0"
$t.Cell(5, 5).Range.Text = "30÷2=15, 0"

$t.Cell(9, 1).Range.Text = "87÷4=21, 3"
$t.Cell(9, 2).Range.Text = "46÷5=9, 1"
$t.Cell(9, 3).Range.Text = "20÷3=6, 2"
$t.Cell(9, 4).Range.Text = "19÷7=2, 5"
$t.Cell(9, 5).Range.Text = "96÷2=48, 0"

$t.Cell(13, 1).Range.Text = "80÷6=13, 2"
$t.Cell(13, 2).Range.Text = "85÷5=17, 0"
$t.Cell(13, 3).Range.Text = "32÷9=3, 5"
$t.Cell(13, 4).Range.Text = "24÷5=4, 4"
$t.Cell(13, 5).Range.Text = "15÷9=1, 6"

$t.Cell(17, 1).Range.Text = "25÷2=12, 1"
$t.Cell(17, 2).Range.Text = "68÷9=7, 5"
$t.Cell(17, 3).Range.Text = "97÷9=10, 7"
$t.Cell(17, 4).Range.Text = "60÷8=7, 4"
$t.Cell(17, 5).Range.Text = "90÷4=22, 2"
